$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.147.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.363.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.78%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.16%  "
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.82%  "
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.356"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.784.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.073.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.351.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "333.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("E20").Value = "  +2.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "62.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "  +1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0741"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("E32").Value = "  +11.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("E38").Value = "  +3.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "145.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "293.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.12%  "
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("E44").Value = "  +2.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.15%  "
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("E47").Value = "  +1.38%  "
$ws.Range("E48").Value = "  +2.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.24%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.31%  "
